$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.404.74'
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').Value = '1.874.45'
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9998'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '238.52'
$ws.Range('E5').Value = '  +0.40%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9997'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4790'
$ws.Range('E7').Value = '  -1.10%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2826'
$ws.Range('E8').Value = '  -2.66%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06523'
$ws.Range('D10').Value = '1.874.59'
$ws.Range('E10').Value = '  -0.85%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07464'
$ws.Range('E11').Value = '  +1.65%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '16.66'
$ws.Range('E12').Value = '  -1.60%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.099'
$ws.Range('E13').Value = '  -1.35%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '88.25'
$ws.Range('E14').Value = '  +0.64%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.6593'
$ws.Range('E15').Value = '  -0.39%  '
$ws.Range('D16').Value = '30.365.29'
$ws.Range('E16').Value = '  +0.17%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '13.34'
$ws.Range('E17').Value = '  -0.59%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.9996'
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007615'
$ws.Range('E19').Value = '  -2.13%  '
$ws.Range('D20').Value = '2.121.25'
$ws.Range('E20').Value = '  -0.64%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.309'
$ws.Range('E21').Value = '  -2.58%  '
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '218.40'
$ws.Range('E23').Value = '  +12.89%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.222'
$ws.Range('E24').Value = '  +0.69%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.344'
$ws.Range('E25').Value = '  -0.45%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '167.49'
$ws.Range('E26').Value = '  +1.90%  '
$ws.Range('E27').Value = '  +1.47%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.977'
$ws.Range('E28').Value = '  +1.88%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.463'
$ws.Range('E29').Value = '  +0.83%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.09435'
$ws.Range('E30').Value = '  +3.29%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.317'
$ws.Range('E31').Value = '  +0.18%  '
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range('E33').Value = '  -0.42%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.208'
$ws.Range('E34').Value = '  +5.73%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.7513'
$ws.Range('E35').Value = '  +2.75%  '
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.01824'
$ws.Range('E37').Value = '  +1.83%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.614'
$ws.Range('E38').Value = '  -1.43%  '
$ws.Range('E39').Value = '  -0.07%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.9059'
$ws.Range('E40').Value = '  -1.45%  '
$ws.Range('E41').Value = '  +1.17%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.902'
$ws.Range('E42').Value = '  +0.14%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.4284'
$ws.Range('E43').Value = '  -0.54%  '
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '7.416'
$ws.Range('E45').Value = '  -0.90%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '64.57'
$ws.Range('E46').Value = '  -0.58%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.1284'
$ws.Range('E47').Value = '  -3.37%  '
$ws.Range('E48').Value = '  -6.79%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '8.969'
$ws.Range('E49').Value = '  -0.07%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '33.59'
$ws.Range('E50').Value = '  -1.05%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.3900'
$ws.Range('E51').Value = '  +1.28%  '
